$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateText($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# ---- Row 5 (05-09-2015 / Saturday) : values only, no other change ----
Set-DateText "A5" "05-09-2015"

# ---- Row 6 (05-10-2015 / Sunday) ----
Set-DateText "A6" "05-10-2015"

# ---- Row 7 (05-11-2015 / Monday) ----
Set-DateText "A7" "05-11-2015"
$ws.Range("C7").Value = "09:17:36"
$ws.Range("D7").Value = "19:06:42"
$ws.Range("E7").Value = 1.0

# ---- Row 8 (05-12-2015 / Tuesday) ----
Set-DateText "A8" "05-12-2015"
$ws.Range("C8").Value = "10:06:52"
$ws.Range("D8").Value = "18:10:57"
$ws.Range("F8").Value = 0.5

# ---- Row 9 (05-13-2015 / Wednesday) ----
Set-DateText "A9" "05-13-2015"
$ws.Range("C9").Value = "08:18:20"
$ws.Range("D9").Value = "18:18:03"
$ws.Range("F9").Value = 0.25

# ---- Row 10 (05-14-2015 / Thursday) : dates only ----
Set-DateText "A10" "05-14-2015"
$ws.Range("C10").Value = "12:22:03"
$ws.Range("D10").Value = "16:04:47"

# ---- Row 11 (05-15-2015 / Friday) -- restyled to style10 ----
Set-DateText "A11" "05-15-2015"
$ws.Range("C11").Value = "14:34:10"
$ws.Range("D11").Value = "16:50:38"
$ws.Range("E11").ClearContents()
$ws.Range("F11").Value = 0.75

# ---- Row 12 (05-16-2015 / Saturday) ----
Set-DateText "A12" "05-16-2015"

# ---- Row 13 (05-17-2015 / Sunday) ----
Set-DateText "A13" "05-17-2015"

# ---- Row 14 (05-18-2015 / Monday) -- restyled to style10 ----
Set-DateText "A14" "05-18-2015"
$ws.Range("C14").Value = "08:19:02"
$ws.Range("D14").Value = "12:57:50"
$ws.Range("E14").ClearContents()

# ---- Row 15 (05-19-2015 / Tuesday) -- restyled to style10 ----
Set-DateText "A15" "05-19-2015"
$ws.Range("C15").Value = "12:39:50"
$ws.Range("D15").Value = "18:14:34"
$ws.Range("E15").ClearContents()
$ws.Range("F15").Value = 0.25

# ---- Row 16 (05-20-2015 / Wednesday) ----
Set-DateText "A16" "05-20-2015"
$ws.Range("C16").Value = "10:54:34"
$ws.Range("D16").Value = "19:03:45"
$ws.Range("F16").ClearContents()

# ---- Row 17 (05-21-2015 / Thursday) -- was a holiday row, now restyled to style10 with real attendance ----
Set-DateText "A17" "05-21-2015"
$ws.Range("C17").Value = "11:29:34"
$ws.Range("D17").Value = "17:37:54"
$ws.Range("F17").Value = 1.0
$ws.Range("P17").ClearContents()

# ---- Row 18 (05-22-2015 / Friday) -- was a holiday row, now restyled to style10 with real attendance ----
Set-DateText "A18" "05-22-2015"
$ws.Range("C18").Value = "10:54:11"
$ws.Range("D18").Value = "13:28:51"
$ws.Range("P18").ClearContents()

# ---- Re-apply the correct direct formatting (cellXfs style) for the rows whose
#      highlight color changed from the diff (orange "style 10" fill, matching
#      the template used by rows 8/10/16). Copy format only, values untouched. ----
$ws.Range("A8:P8").Copy()
$ws.Range("A11:P11").PasteSpecial(-4122)
$ws.Range("A14:P14").PasteSpecial(-4122)
$ws.Range("A15:P15").PasteSpecial(-4122)
$ws.Range("A17:P17").PasteSpecial(-4122)
$ws.Range("A18:P18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 19/20: count & sum ranges now extend through row 18 ----
$ws.Range("E19").Formula = "=COUNT(E5:E18)"
$ws.Range("E20").Formula = "=SUM(E5:E18)"

# ---- Row 22: H/I sums now extend through row 18; I22 also picks up the
#      "style 17" direct formatting already used by H22 ----
$ws.Range("H22").Formula = "=SUM(H5:H18)"
$ws.Range("I22").Formula = "=SUM(I5:I18)"
$ws.Range("H22").Copy()
$ws.Range("I22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 23: literal total-absences value corrected ----
$ws.Range("I23").Value = 4.0
